$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains weekly price records for "Perejil" (Hortaliza), two rows
# (Primera / Segunda) per reporting date, sorted with the newest "window" of
# dates starting at row 14 (rows 2-13 hold three fixed historical weeks that
# are not touched). This edit adds a new, more recent weekly record
# (2023-06-15) at the top of that rotating window (rows 14-15), which pushes
# every subsequent pair of rows down by one pair (2 rows); the oldest pair
# that used to sit at rows 224-225 ends up at the newly created rows 226-227.

# Insert two blank rows at row 14, shifting rows 14-225 down to 16-227.
$ws.Range("A14:R15").EntireRow.Insert()

# The old content that used to live at rows 14-15 is now at rows 16-17.
# Duplicate it back into the new rows 14-15 as a starting point.
$src = $ws.Range("A16:R17")
$dst = $ws.Range("A14:R15")
$src.Copy($dst)

# Update only the date of the newly inserted pair to the new reporting date.
$ws.Range("D14").Value = [DateTime]"2023-06-15"
$ws.Range("D15").Value = [DateTime]"2023-06-15"
